# "Se implemento cola dinamica" - add the new task row (row 8) to the
# Metricas sheet: a new shared-string label "implementar Cola Dinamica"
# plus its estimated/actual line counts and start/end times. The
# dependent totals (row 11), summary (rows 12, 15-20) and the pie chart
# that reads C19:C20 are formula-driven and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")
$ws.Activate()

$ws.Range("A8").Value = "implementar Cola Dinamica"
$ws.Range("B8").Value = 40
$ws.Range("C8").Value = 36
$ws.Range("D8").Value = 0.020833333333333332
$ws.Range("E8").Value = 0.11944444444444445
$ws.Range("F8").Value = 0.1451388888888889
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

# Recalculate so every dependent formula (totals, summary %, chart source
# cells) carries a fresh cached value before the workbook is written out.
$excel.CalculateFull()

# The author's selection moved to C14 before saving.
$ws.Range("C14").Select()

$wb.Save()
